# Insert a new blank row above row 7 (pushes "Test on Adam optimization"
# and everything below it down by one row), then move the selection to
# A37, matching the author's final cursor position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert() | Out-Null

$ws.Range("A37").Select() | Out-Null
